$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'66.658.89"
$ws.Range('E2').Value = "'  -0.43%  "

$ws.Range('D3').Value = "'3.068.25"
$ws.Range('E3').Value = "'  -1.54%  "

$ws.Range('E4').Value = "'  -0.06%  "

$ws.Range('D5').Value = "'574.76"
$ws.Range('E5').Value = "'  -0.53%  "

$ws.Range('D6').Value = "'169.53"
$ws.Range('E6').Value = "'  -1.10%  "

$ws.Range('D7').Value = "'0.998"
$ws.Range('E7').Value = "'  -0.14%  "

$ws.Range('D8').Value = "'3.067.24"
$ws.Range('E8').Value = "'  -1.48%  "

$ws.Range('D9').Value = "'0.511"
$ws.Range('E9').Value = "'  -2.02%  "

$ws.Range('D10').Value = "'6.36"
$ws.Range('E10').Value = "'  -1.68%  "

$ws.Range('E11').Value = "'  -1.75%  "

$ws.Range('D12').Value = "'0.467"
$ws.Range('E12').Value = "'  -3.23%  "

$ws.Range('E13').Value = "'  -2.62%  "

$ws.Range('D14').Value = "'35.69"
$ws.Range('E14').Value = "'  -4.04%  "

$ws.Range('E15').Value = "'  -1.85%  "

$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = "'3.577.37"
$ws.Range('E16').Value = "'  -1.51%  "

$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = "'66.623.54"
$ws.Range('E17').Value = "'  -0.44%  "

$ws.Range('D18').Value = "'6.99"
$ws.Range('E18').Value = "'  -2.26%  "

$ws.Range('D19').Value = "'16.83"
$ws.Range('E19').Value = "'  +2.95%  "

$ws.Range('D20').Value = "'3.064.21"
$ws.Range('E20').Value = "'  -1.64%  "

$ws.Range('D21').Value = "'492.32"
$ws.Range('E21').Value = "'  +3.25%  "

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'7.72"
$ws.Range('E22').Value = "'  -3.21%  "

$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').Value = "'0.688"
$ws.Range('E23').Value = "'  -3.62%  "

$ws.Range('D24').Value = "'82.78"
$ws.Range('E24').Value = "'  -1.50%  "

$ws.Range('D25').Value = "'12.67"
$ws.Range('E25').Value = "'  -5.58%  "

$ws.Range('D26').Value = "'2.20"
$ws.Range('E26').Value = "'  -4.11%  "

$ws.Range('D27').Value = "'10.16"
$ws.Range('E27').Value = "'  +1.05%  "

$ws.Range('E28').Value = "'  -0.02%  "

$ws.Range('D29').Value = "'7.76"
$ws.Range('E29').Value = "'  -1.78%  "

$ws.Range('D30').Value = "'2.27"
$ws.Range('E30').Value = "'  -5.18%  "

$ws.Range('E31').Value = "'  -2.55%  "

$ws.Range('D32').Value = "'27.52"
$ws.Range('E32').Value = "'  -3.60%  "

$ws.Range('D33').Value = "'0.112"
$ws.Range('E33').Value = "'  -3.17%  "

$ws.Range('D34').Value = "'0.0₃0910"
$ws.Range('E34').Value = "'  -3.16%  "

$ws.Range('D35').Value = "'0.999"
$ws.Range('E35').Value = "'  -0.03%  "

$ws.Range('D36').Value = "'0.950"
$ws.Range('E36').Value = "'  -2.47%  "

$ws.Range('D37').Value = "'5.58"
$ws.Range('E37').Value = "'  -4.95%  "

$ws.Range('D38').Value = "'46.43"
$ws.Range('E38').Value = "'  -1.37%  "

$ws.Range('E39').Value = "'  +0.02%  "

$ws.Range('E40').Value = "'  -5.57%  "

$ws.Range('E41').Value = "'  -3.39%  "

$ws.Range('D42').Value = "'8.31"
$ws.Range('E42').Value = "'  -4.65%  "

$ws.Range('D43').Value = "'2.755.41"
$ws.Range('E43').Value = "'  -3.22%  "

$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = "'0.0345"
$ws.Range('E44').Value = "'  -3.43%  "

$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').Value = "'135.67"
$ws.Range('E45').Value = "'  -0.26%  "

$ws.Range('D46').Value = "'366.47"
$ws.Range('E46').Value = "'  -5.04%  "

$ws.Range('D47').Value = "'2.48"
$ws.Range('E47').Value = "'  -4.22%  "

$ws.Range('E48').Value = "'  +0.00%  "

$ws.Range('D49').Value = "'24.53"
$ws.Range('E49').Value = "'  -1.19%  "

$ws.Range('E50').Value = "'  -2.20%  "

$ws.Range('E51').Value = "'  -1.94%  "
